# Auto-generated edit script applying the Mateus_Profits market-data refresh.
# Each sheet corresponds to a crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For every touched row we rewrite the refreshed market columns (H:N); a couple of
# rows gain or lose a profit cell (M/N) depending on whether HQ/NQ pricing applies.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 15000
$ws.Range("I18").Value = 15000
$ws.Range("K18").Value = 15000
$ws.Range("M18").Value = -14716
$ws.Range("H38").Value = 1286.4445
$ws.Range("J38").Value = 3499.6667
$ws.Range("L38").Value = 10499.0001
$ws.Range("N38").Value = -11243.0001
$ws.Range("H41").Value = 147.92308
$ws.Range("J41").Value = 59.5
$ws.Range("L41").Value = 59.5
$ws.Range("N41").Value = -939.5
$ws.Range("H43").Value = 6123.6665
$ws.Range("J43").Value = 7692.5
$ws.Range("L43").Value = 7692.5
$ws.Range("N43").Value = -7830.5
$ws.Range("H58").Value = 54
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H62").Value = 19950.666
$ws.Range("I62").Value = 18079.572
$ws.Range("K62").Value = 18079.572
$ws.Range("M62").Value = -17455.572
$ws.Range("H65").Value = 19950.666
$ws.Range("I65").Value = 18079.572
$ws.Range("K65").Value = 90397.86
$ws.Range("M65").Value = -87277.86
$ws.Range("H132").Value = 1918.8636
$ws.Range("I132").Value = 1836.1052
$ws.Range("K132").Value = 5508.3156
$ws.Range("M132").Value = -2978.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 30000
$ws.Range("J58").Value = 30000
$ws.Range("L58").Value = 30000
$ws.Range("N58").Value = -30860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 35000
$ws.Range("I82").Value = 20000
$ws.Range("K82").Value = 20000
$ws.Range("M82").Value = -19617
$ws.Range("H85").Value = 35000
$ws.Range("I85").Value = 20000
$ws.Range("K85").Value = 20000
$ws.Range("M85").Value = -18674

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1313.1818
$ws.Range("J5").Value = 4333.3335
$ws.Range("L5").Value = 4333.3335
$ws.Range("N5").Value = -4557.3335
$ws.Range("H16").Value = 3198.9473
$ws.Range("I16").Value = 2487.9167
$ws.Range("K16").Value = 2487.9167
$ws.Range("M16").Value = -2200.9167
$ws.Range("H25").Value = 9997
$ws.Range("I25").Value = 9997
$ws.Range("K25").Value = 9997
$ws.Range("M25").Value = -9823
$ws.Range("H31").Value = 5645.273
$ws.Range("J31").Value = 6914.6665
$ws.Range("L31").Value = 6914.6665
$ws.Range("N31").Value = -7504.6665
$ws.Range("H34").Value = 5645.273
$ws.Range("J34").Value = 6914.6665
$ws.Range("L34").Value = 6914.6665
$ws.Range("N34").Value = -7318.6665
$ws.Range("H41").Value = 20785
$ws.Range("I41").Value = 24199
$ws.Range("K41").Value = 24199
$ws.Range("M41").Value = -23771
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H60").Value = 27995
$ws.Range("I60").Value = 27985
$ws.Range("J60").Value = 28000
$ws.Range("K60").Value = 27985
$ws.Range("L60").Value = 28000
$ws.Range("M60").Value = -27474
$ws.Range("N60").Value = -29022
$ws.Range("H74").Value = 37636.668
$ws.Range("J74").Value = 37636.668
$ws.Range("L74").Value = 37636.668
$ws.Range("N74").Value = -39384.668
$ws.Range("H77").Value = 37636.668
$ws.Range("J77").Value = 37636.668
$ws.Range("L77").Value = 112910.004
$ws.Range("N77").Value = -121646.004
$ws.Range("H86").Value = 8569.714
$ws.Range("I86").Value = 7497.25
$ws.Range("J86").Value = 9999.666999999999
$ws.Range("K86").Value = 7497.25
$ws.Range("L86").Value = 9999.666999999999
$ws.Range("M86").Value = -6374.25
$ws.Range("N86").Value = -12245.667
$ws.Range("H89").Value = 8569.714
$ws.Range("I89").Value = 7497.25
$ws.Range("J89").Value = 9999.666999999999
$ws.Range("K89").Value = 37486.25
$ws.Range("L89").Value = 49998.335
$ws.Range("M89").Value = -31870.25
$ws.Range("N89").Value = -61230.335
$ws.Range("H105").Value = 2367.5
$ws.Range("I105").Value = 2531.4285
$ws.Range("K105").Value = 2531.4285
$ws.Range("M105").Value = -784.4285
$ws.Range("H107").Value = 474
$ws.Range("J107").Value = 592.2
$ws.Range("L107").Value = 592.2
$ws.Range("N107").Value = -4432.2
$ws.Range("H111").Value = 68613
$ws.Range("J111").Value = 68613
$ws.Range("L111").Value = 68613
$ws.Range("N111").Value = -76793
$ws.Range("H113").Value = 3198.9473
$ws.Range("I113").Value = 2487.9167
$ws.Range("K113").Value = 2487.9167
$ws.Range("M113").Value = -317.9167000000002
$ws.Range("H122").Value = 3582.4443
$ws.Range("I122").Value = 3352.4443
$ws.Range("K122").Value = 10057.3329
$ws.Range("M122").Value = -7607.332900000001
$ws.Range("H134").Value = 5317.636
$ws.Range("I134").Value = 4360.1665
$ws.Range("K134").Value = 13080.4995
$ws.Range("M134").Value = -10545.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16173574
$ws.Range("I4").Value = 14429059
$ws.Range("K4").Value = 43287177
$ws.Range("M4").Value = -43287065

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 35000
$ws.Range("J26").Value = 35000
$ws.Range("L26").Value = 35000
$ws.Range("N26").Value = -35560
$ws.Range("H43").Value = 15000
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 20000
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = -4849
$ws.Range("N43").Value = -20302
$ws.Range("H50").Value = 35000
$ws.Range("J50").Value = 35000
$ws.Range("L50").Value = 35000
$ws.Range("N50").Value = -35996
$ws.Range("H70").Value = 13128.571
$ws.Range("I70").Value = 16000
$ws.Range("J70").Value = 12650
$ws.Range("K70").Value = 16000
$ws.Range("L70").Value = 12650
$ws.Range("M70").Value = -15730
$ws.Range("N70").Value = -13190
$ws.Range("H73").Value = 13128.571
$ws.Range("I73").Value = 16000
$ws.Range("J73").Value = 12650
$ws.Range("K73").Value = 16000
$ws.Range("L73").Value = 12650
$ws.Range("M73").Value = -15064
$ws.Range("N73").Value = -14522
$ws.Range("H80").Value = 4489.75
$ws.Range("I80").Value = 4223.4546
$ws.Range("J80").Value = 5075.6
$ws.Range("K80").Value = 4223.4546
$ws.Range("L80").Value = 5075.6
$ws.Range("M80").Value = -3225.4546
$ws.Range("N80").Value = -7071.6
$ws.Range("H83").Value = 4489.75
$ws.Range("I83").Value = 4223.4546
$ws.Range("J83").Value = 5075.6
$ws.Range("K83").Value = 21117.273
$ws.Range("L83").Value = 25378
$ws.Range("M83").Value = -16125.273
$ws.Range("N83").Value = -35362
$ws.Range("H113").Value = 4366.8
$ws.Range("I113").Value = 3499.0833
$ws.Range("K113").Value = 3499.0833
$ws.Range("M113").Value = -1329.0833
$ws.Range("H122").Value = 2898.9546
$ws.Range("I122").Value = 2688.85
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 8066.549999999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -5616.549999999999
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 4603
$ws.Range("J126").Value = 6089.8
$ws.Range("L126").Value = 18269.4
$ws.Range("N126").Value = -23209.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4481.125
$ws.Range("I22").Value = 3370
$ws.Range("K22").Value = 3370
$ws.Range("M22").Value = -3075
$ws.Range("H27").Value = 4481.125
$ws.Range("I27").Value = 3370
$ws.Range("K27").Value = 3370
$ws.Range("M27").Value = -3263
$ws.Range("H40").Value = 3154.7778
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3272
$ws.Range("H122").Value = 2585
$ws.Range("I122").Value = 2590.25
$ws.Range("K122").Value = 7770.75
$ws.Range("M122").Value = -5320.75
$ws.Range("H136").Value = 3650.5908
$ws.Range("I136").Value = 3121.353
$ws.Range("J136").Value = 5450
$ws.Range("K136").Value = 9364.059000000001
$ws.Range("L136").Value = 16350
$ws.Range("M136").Value = -6814.059000000001
$ws.Range("N136").Value = -21450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1666.5555
$ws.Range("H107").Value = 1388.6316
$ws.Range("I107").Value = 711.7692
$ws.Range("J107").Value = 2855.1667
$ws.Range("K107").Value = 2135.3076
$ws.Range("L107").Value = 8565.500100000001
$ws.Range("M107").Value = -215.3076000000001
$ws.Range("N107").Value = -12405.5001
